# Realestate Update resale numbers 2024-01-10 17:40
# Appends a new data row (row 43) to the CityResaleNum sheet with the
# latest resale numbers snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

# Date/time/week-looking text must be kept as text (matches how the
# existing rows store these columns as plain strings, not Excel dates
# or numbers) - format the cells as Text *before* assigning the value so
# Excel does not auto-convert them to date/time serials or numbers, then
# drop the temporary formatting again so the new row keeps the sheet's
# default (unstyled) look, same as every other data row.
$leadCells = $ws.Range("A$($row):D$($row)")
$leadCells.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-10"
$ws.Cells.Item($row, 2).Value = "17:40:08"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "01"

$leadCells.ClearFormats()

$ws.Cells.Item($row, 5).Value = 139614
$ws.Cells.Item($row, 6).Value = 142642
$ws.Cells.Item($row, 7).Value = 172105
$ws.Cells.Item($row, 8).Value = 147993
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119131
$ws.Cells.Item($row, 11).Value = 224752
$ws.Cells.Item($row, 12).Value = 251281
$ws.Cells.Item($row, 13).Value = 185036
$ws.Cells.Item($row, 14).Value = 110490
$ws.Cells.Item($row, 15).Value = 40754
$ws.Cells.Item($row, 16).Value = 30861
$ws.Cells.Item($row, 17).Value = 72754
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42114
$ws.Cells.Item($row, 20).Value = -1
